# Update the division problems in the practice-sheet table.
# Each populated row of the table (rows 1, 5, 9, 13, 17) holds five
# "two-digit ÷ one-digit =" expressions, one per column. We address each
# cell directly by (row, column) so there is no ambiguity even though some
# new values coincide with other cells' old values (e.g. "23÷6=").

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="83÷4="},
    @{Row=1;  Col=2; Text="56÷7="},
    @{Row=1;  Col=3; Text="33÷5="},
    @{Row=1;  Col=4; Text="79÷8="},
    @{Row=1;  Col=5; Text="94÷7="},

    @{Row=5;  Col=1; Text="32÷3="},
    @{Row=5;  Col=2; Text="18÷7="},
    @{Row=5;  Col=3; Text="54÷5="},
    @{Row=5;  Col=4; Text="18÷2="},
    @{Row=5;  Col=5; Text="22÷6="},

    @{Row=9;  Col=1; Text="40÷8="},
    @{Row=9;  Col=2; Text="70÷5="},
    @{Row=9;  Col=3; Text="20÷3="},
    @{Row=9;  Col=4; Text="23÷6="},
    @{Row=9;  Col=5; Text="29÷6="},

    @{Row=13; Col=1; Text="60÷9="},
    @{Row=13; Col=2; Text="84÷7="},
    @{Row=13; Col=3; Text="11÷2="},
    @{Row=13; Col=4; Text="66÷5="},
    @{Row=13; Col=5; Text="33÷6="},

    @{Row=17; Col=1; Text="90÷7="},
    @{Row=17; Col=2; Text="74÷3="},
    @{Row=17; Col=3; Text="79÷5="},
    @{Row=17; Col=4; Text="40÷2="},
    @{Row=17; Col=5; Text="52÷4="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
